$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Plain "want to go" count bumps (column F) on existing rows.
$ws1Updates = @{
    2  = 825
    3  = 564
    4  = 312
    6  = 1158
    7  = 339
    8  = 54
    11 = 1210
    14 = 901
    15 = 884
    17 = 75
    20 = 790
    21 = 1746
    22 = 3012
    23 = 881
    25 = 2257
    26 = 670
    28 = 3109
    29 = 623
    30 = 575
    33 = 740
}
foreach ($r in $ws1Updates.Keys) {
    $ws1.Cells.Item($r, 6).Value2 = $ws1Updates[$r]
}

# Insert a new event row at row 36 (shifts the old rows 36-45 down to 37-46).
$ws1.Rows.Item(36).Insert(-4121)

# Copy the row-index column's formatting onto the freshly inserted cell so it
# keeps the bordered/bold/centered look used throughout column A.
$ws1.Range("A37").Copy()
$ws1.Range("A36").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Cells.Item(36, 2).Value2 = "'2024-05-26"
$ws1.Cells.Item(36, 3).Value2 = "杭州·恋与深空×恋与制作人only"
$ws1.Cells.Item(36, 4).Value2 = "望江东路333号 杭州瑞莱克斯大酒店"
$ws1.Cells.Item(36, 5).Value2 = "2024.05.26 10:00-05.26 17:00"
$ws1.Cells.Item(36, 6).Value2 = 5
$ws1.Cells.Item(36, 7).Value2 = 60
$ws1.Cells.Item(36, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=84077"
$ws1.Cells.Item(36, 9).Value2 = "//i1.hdslb.com/bfs/openplatform/202404/V6V4Pppv1712736555042.jpeg"

# Column A is a plain sequential row index (0-based); keep it sequential for
# the inserted row through the end of the sheet.
for ($r = 36; $r -le 46; $r++) {
    $ws1.Cells.Item($r, 1).Value2 = $r - 1
}

# The rows that got pushed down by the insert also received their own
# "want to go" count refresh (new row numbers, after the shift).
$ws1ShiftedUpdates = @{
    38 = 1106
    39 = 1794
    40 = 401
    42 = 560
    43 = 200
    44 = 134
    45 = 183
    46 = 47
}
foreach ($r in $ws1ShiftedUpdates.Keys) {
    $ws1.Cells.Item($r, 6).Value2 = $ws1ShiftedUpdates[$r]
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(12, 6).Value2 = 84

# Insert a new event row at row 15 (shifts old rows 15-16 down to 16-17).
$ws2.Rows.Item(15).Insert(-4121)

$ws2.Range("A16").Copy()
$ws2.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Cells.Item(15, 2).Value2 = "'2024-06-23"
$ws2.Cells.Item(15, 3).Value2 = "杭州·《亚米·跨越二次元》ACG经典动漫视听音乐会"
$ws2.Cells.Item(15, 4).Value2 = "金沙大道681号 金沙湖大剧院"
$ws2.Cells.Item(15, 5).Value2 = "2024.06.23 19:30-06.23 21:10"
$ws2.Cells.Item(15, 6).Value2 = 0
$ws2.Cells.Item(15, 7).Value2 = 80
$ws2.Cells.Item(15, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=84041"
$ws2.Cells.Item(15, 9).Value2 = "//i2.hdslb.com/bfs/openplatform/202404/UhUuHfad1712564787267.jpeg"

for ($r = 15; $r -le 17; $r++) {
    $ws2.Cells.Item($r, 1).Value2 = $r - 1
}

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types combined) - value-only refresh, no new rows.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4Updates = @{
    2  = 825
    3  = 564
    4  = 312
    6  = 1158
    7  = 339
    8  = 54
    10 = 1210
    12 = 901
    13 = 884
    16 = 75
    19 = 790
    20 = 1746
    21 = 3012
    22 = 881
    25 = 2257
    26 = 3109
    27 = 623
    28 = 576
    35 = 84
    36 = 740
    41 = 1106
    42 = 1794
    44 = 401
    45 = 560
    46 = 200
    47 = 134
    48 = 183
    49 = 47
}
foreach ($r in $ws4Updates.Keys) {
    $ws4.Cells.Item($r, 6).Value2 = $ws4Updates[$r]
}
